$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newCourses = @("MATH094", "MATH095", "MATH096", "MATH098", "ENGL097", "ENGL098")

$startRow = 514
for ($i = 0; $i -lt $newCourses.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newCourses[$i]
}

$ws.Range("A520").Select()
$excel.ActiveWindow.ScrollRow = 502
